# Fixed naive component forecaster bug - Presentation state 11.02.
#
# A new evaluation quarter (Q0) was computed upstream, which pushes the
# existing rolling-window rows of error metrics (ME, MAE, MSE, RMSE, SE)
# down by one position and bumps up every "N" (sample size) count by one
# - the oldest row (old Q9 row) drops off the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Shift the existing metric rows (B2:G10) down into (B3:G11).
#    This reproduces old row r becoming new row r+1 for every metric
#    column (B..F) as well as the N column (G), row labels in column A
#    stay put (Q0..Q9 labels are unaffected by the data refresh).
$src = $ws.Range("B2:G10")
$dst = $ws.Range("B3:G11")
$dst.Value2 = $src.Value2

# 2) Write the freshly computed values for the newest row (row 2), and
#    bump its sample size (N) to 15.
$ws.Range("B2").Value2 = 0.2015370511150554
$ws.Range("C2").Value2 = 0.3501553535809984
$ws.Range("D2").Value2 = 0.2617601871928103
$ws.Range("E2").Value2 = 0.5116250455097076
$ws.Range("F2").Value2 = 0.4867637343656181
$ws.Range("G2").Value2 = 15
